# Update hashcode values in column B for the specified rows (sharedStrings-backed cells)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("hashcode.csv")

$updates = @{
    89 = "3a425473b901d99eeb2f8f05d1a7a9da"
    99 = "7332e19db9d80de1248db805e60f9312"
    110 = "a0cab0e46f110ea81f706b2fc5953f20"
    154 = "dadb7be999dbd63f806299bfafbc6261"
    160 = "25264021f32130c246ff1dcdeec483d0"
    169 = "4da83de0fa8baa0c3e34ef948fa497bf"
    222 = "b2c2d7b0c6e1e482e2baebfaa3e80238"
    227 = "811e4b110a2cffba77fce045c7017d73"
    229 = "67e8de9238b1d980854c534789e8446c"
    232 = "869c621bbced2dd1e9009bcaac137d49"
    281 = "181895aa68478a8ce5e37e3a6123fdf6"
    284 = "afc91a4d0896544a39504d970bebe301"
    338 = "c16252edd9bbad81bece7e1e437aeca5"
    423 = "0841f66eec1f7caf51680bed6f5054c6"
    486 = "7c7e26fef28b133513b0e1d817db11ed"
    511 = "b3c0471f6ab03fe79ed3515cd46b22cc"
    526 = "46abcc7d85f2732d753478da077c6dad"
    545 = "6872b106d46507f66af37d33523f76f9"
    559 = "a43aad2a42277be6fc85233bafe81f21"
    565 = "2ba2af195a7150411e9edbf214040e44"
    578 = "c2773ef09b571a4d55e3f514b1138e7d"
    584 = "90e9978e5fac4cdc1c413f6cc4049a3c"
    596 = "db79560a07b943a028661bf9ac58f8cf"
    677 = "16b63d480f3d50d78a869c19ab998727"
    692 = "87f7d8c8d5f14748512c9245c79f6ea6"
    697 = "e992428de39ad6cc52cb72f089587295"
    712 = "c73244e4d02da93b2f5418460dd36c9d"
    715 = "d174fa8fbca0c777f41402c2571309ad"
    780 = "7b32c2e2138ad20d6de90800ca768f42"
    823 = "1240d1925d5bb6781d888325f1408e49"
    827 = "18959c8132fbe58132b63e2ed262ede7"
    828 = "683ad9d5a62eedccab952d06bed5a4f7"
    837 = "c23d1d2e9e89bd032e026d27dfcc8827"
    839 = "97010d418992034607b9ffb8ac4a8020"
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 2).Value = $updates[$row]
}

Write-Output "Updated $($updates.Count) hashcode cells."
